$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for the market. Insert a new row at
# row 14 (pushing the existing rows 14-27 down to 15-28) and populate it
# with the new record's data.
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C14").Value = "Arica y Parinacota"
$ws.Range("D14").Value = 44580
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = 100112028
$ws.Range("G14").Value = "Sandia"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 1200
$ws.Range("K14").Value = 380
$ws.Range("L14").Value = 400
$ws.Range("M14").Value = 390
$ws.Range("N14").Value = "$/kilo (volumen en unidades)"
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 390
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"
